$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update matchup averages for spring 24 regular season complete
$ws.Range("F3").Value = 1.22
$ws.Range("G4").Value = 1.03
$ws.Range("D5").Value = 1.34
$ws.Range("G6").Value = 0.97
$ws.Range("D7").Value = 1.68
$ws.Range("F7").Value = 1.5
